$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Column width changes ---
# Raw OOXML `width` = COM `ColumnWidth` + 0.8333333333333334 (5/6 char), so
# back the offset out to land on the exact target raw widths (11 and 9.5).
$ws.Columns.Item(5).ColumnWidth = 11 - 0.8333333333333334
$ws.Columns.Item(13).ColumnWidth = 9.5 - 0.8333333333333334

# --- Add new ratio columns / values to rows 5, 7, 9, 12, 13 ---

# Row 5
$ws.Range("C5").Value = 55105.8
$ws.Range("D5").Formula = "=H5/C5"
$ws.Range("E5").Value = 129011106
$ws.Range("F5").Formula = "=J5/E5"
$ws.Range("I5").Formula = "=H5/H5"
$ws.Range("K5").Formula = "=J5/J5"
$ws.Range("M5").Value = 102864
$ws.Range("N5").Formula = "=H5/M5"
$ws.Range("P5").Formula = "=J5/O5"

# Row 7
$ws.Range("D7").Formula = "=H7/C7"
$ws.Range("F7").Formula = "=J7/E7"
$ws.Range("I7").Formula = "=H7/H7"
$ws.Range("K7").Formula = "=J7/J7"
$ws.Range("N7").Formula = "=H7/M7"
$ws.Range("P7").Formula = "=J7/O7"

# Row 9
$ws.Range("D9").Formula = "=H9/C9"
$ws.Range("F9").Formula = "=J9/E9"
$ws.Range("I9").Formula = "=H9/H9"
$ws.Range("K9").Formula = "=J9/J9"
$ws.Range("N9").Formula = "=H9/M9"
$ws.Range("P9").Formula = "=J9/O9"

# Row 12
$ws.Range("D12").Formula = "=H12/C12"
$ws.Range("F12").Formula = "=J12/E12"
$ws.Range("I12").Formula = "=H12/H12"
$ws.Range("K12").Formula = "=J12/J12"
$ws.Range("N12").Formula = "=H12/M12"
$ws.Range("P12").Formula = "=J12/O12"

# Row 13
$ws.Range("D13").Formula = "=H13/C13"
$ws.Range("F13").Formula = "=J13/E13"
$ws.Range("I13").Formula = "=H13/H13"
$ws.Range("K13").Formula = "=J13/J13"
$ws.Range("N13").Formula = "=H13/M13"
$ws.Range("P13").Formula = "=J13/O13"

# --- Insert 3 new blank rows before row 16 (pushes old 16.. down by 3) ---
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# New rows 16, 17, 18 should look like the blank separator row (A col style=5, B/G/L/Q col style=1)
$ws.Range("A12").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)

$ws.Range("B12").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)

$ws.Range("G12").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G18").PasteSpecial(-4122)

$ws.Range("L12").Copy()
$ws.Range("L16").PasteSpecial(-4122)
$ws.Range("L17").PasteSpecial(-4122)
$ws.Range("L18").PasteSpecial(-4122)

$ws.Range("Q12").Copy()
$ws.Range("Q16").PasteSpecial(-4122)
$ws.Range("Q17").PasteSpecial(-4122)
$ws.Range("Q18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- New note row 31 with new shared string ---
$ws.Range("B31").Value = "*stream would not run on 10 locales"

# --- Selection / view state ---
$ws.Range("B34").Select()

# Best-effort: nudge the saved window x-position to match the source edit
# (xWindow -140 -> -160). Harmless no-op if the host doesn't expose a
# Window object backed by the OOXML bookViews/workbookView part.
try {
    $excel.ActiveWindow.Left = -160
} catch {}
